# Update the "as_of_utc" timestamp column (AA, rows 2-26) on the
# "Главные" and "Линейные" sheets from 2025-11-18 07:06:23 to
# 2025-11-18 10:36:15.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-18 10:36:15"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("AA2:AA26").Value = $newTimestamp
}
